$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + date range) ---
$ws.Range("A8").Value = "Volume 32   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"

# --- Crime-data grid updates ---
$ws.Range("M15").Value = -31.25
$ws.Range("N15").Value = -42.105263157894
$ws.Range("C16").Value = 2
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -44.444444444444
$ws.Range("I16").Value = 80
$ws.Range("J16").Value = 97
$ws.Range("K16").Value = -17.525773195876
$ws.Range("L16").Value = -12.087912087912
$ws.Range("M16").Value = -55.056179775280
$ws.Range("N16").Value = -85.480943738657
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 3
$ws.Range("D17").NumberFormat = "#,##0"
$ws.Range("E17").Value = -66.666666666666
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 23.076923076923
$ws.Range("I17").Value = 125
$ws.Range("J17").Value = 147
$ws.Range("K17").Value = -14.965986394557
$ws.Range("L17").Value = -7.407407407407
$ws.Range("M17").Value = 22.549019607843
$ws.Range("N17").Value = -51.923076923076
$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 1
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("E18").Value = 0
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 59
$ws.Range("K18").Value = -10.169491525423
$ws.Range("L18").Value = -32.051282051282
$ws.Range("M18").Value = -76.855895196506
$ws.Range("N18").Value = -94.301075268817
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 2.325581395348
$ws.Range("I19").Value = 417
$ws.Range("J19").Value = 471
$ws.Range("K19").Value = -11.464968152866
$ws.Range("L19").Value = -20.872865275142
$ws.Range("M19").Value = 5.037783375314
$ws.Range("N19").Value = -20.419847328244
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -41.176470588235
$ws.Range("J20").Value = 145
$ws.Range("K20").Value = -13.793103448275
$ws.Range("L20").Value = 23.762376237623
$ws.Range("M20").Value = -4.580152671755
$ws.Range("N20").Value = -94.343891402714
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -35.294117647058
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -8.988764044943
$ws.Range("I21").Value = 814
$ws.Range("J21").Value = 930
$ws.Range("K21").Value = -12.473118279569
$ws.Range("L21").Value = -13.953488372093
$ws.Range("M21").Value = -22.770398481973
$ws.Range("N21").Value = -81.943212067435
$ws.Range("D23").Value = 2
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G23").Value = 2
$ws.Range("G23").NumberFormat = "#,##0"
$ws.Range("H23").Value = -100
$ws.Range("H23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = -42.307692307692
$ws.Range("L23").Value = -34.782608695652
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -10.344827586206
$ws.Range("F24").Value = 87
$ws.Range("H24").Value = -43.137254901960
$ws.Range("I24").Value = 1014
$ws.Range("J24").Value = 1318
$ws.Range("K24").Value = -23.065250379362
$ws.Range("L24").Value = -2.312138728323
$ws.Range("M24").Value = 26.591760299625
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = -13.636363636363
$ws.Range("F25").Value = 66
$ws.Range("G25").Value = 133
$ws.Range("H25").Value = -50.375939849624
$ws.Range("I25").Value = 718
$ws.Range("J25").Value = 1096
$ws.Range("K25").Value = -34.489051094890
$ws.Range("L25").Value = -3.753351206434
$ws.Range("C26").Value = 6
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = 36.842105263157
$ws.Range("I26").Value = 251
$ws.Range("J26").Value = 251
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 17.289719626168
$ws.Range("M26").Value = -11.929824561403
$ws.Range("L28").Value = -12.5
$ws.Range("G31").Value = 1
